$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell is stored as literal text (matching the source data feed's
# inline-string format), so values are entered with a leading quote to
# prevent Excel from auto-converting numeric-looking / percent-looking
# text into real numbers, then the quote-prefix style is cleared so the
# cell format stays "General" / default, matching the original cells.

$ws.Range("D2").Value = "`'291.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "`'-6.14%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "`'39.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "`'-2.58%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "`'5.020"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "`'-3.53%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "`'0.07365"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "`'-4.11%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "`'4.283"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "`'-0.26%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "`'1.553"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "`'-10.86%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "`'0.9121"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "`'-1.95%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "`'0.1198"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "`'-6.42%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "`'0.1749"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "`'-3.70%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "`'0.08708"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "`'-4.58%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "`'0.04163"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "`'-0.09%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "`'0.03%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "`'0.001272"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "`'-1.05%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "`'0.005857"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "`'-0.36%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "`'3.385"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "`'0.96%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D18").Value = "`'0.3299"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "`'-0.65%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "`'7.538"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "`'0.1351"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "`'0.01%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "`'0.2882"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "`'6.10%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "`'0.03829"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "`'-4.70%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "`'0.001276"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "`'0.81%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "`'0.003686"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "`'-10.05%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "`'0.0001281"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "`'0.81%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "`'0.0003726"
$ws.Range("D26").Style = "Normal"
$ws.Range("D38").Value = "`'0.02331"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "`'-7.96%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "`'0.05016"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "`'-5.60%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "`'0.007671"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "`'-2.37%"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "`'148.86%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "`'0.1274"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "`'-3.01%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "`'0.007393"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "`'11.20%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "`'0.007678"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "`'-5.13%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "`'0.3160"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "`'2.48%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "`'0.00006507"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "`'-3.86%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "`'0.03%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "`'12.03%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "`'0.004204"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "`'35.46%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "`'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "`'0.03%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "`'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "`'0.03%"
$ws.Range("E51").Style = "Normal"

Write-Host "Updated 70 cells"
